$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 71.7556381225586
$ws.Range("C2").Value = 7.0
$ws.Range("D2").Value = 40.55263137817383
$ws.Range("E2").Value = 57.85714340209961
$ws.Range("H2").Value = 6.25
